$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$tl = $s.TimeLine
$ms = $tl.MainSequence
Write-Output ("MainSequence count: " + $ms.Count)
